$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A, shifting the existing data (old A:D, the
# four BERTScore stat columns) one column to the right, to B:E. Column A's
# header cell (row 1) is left blank, matching the new layout.
$ws.Columns("A:A").Insert(-4161)

# Force column A's data rows to be treated as text so labels such as
# "25%" aren't auto-converted to numeric percentages.
$ws.Range("A2:A9").NumberFormat = "@"

# Fill in the pandas describe()-style index labels in column A, rows 2-9.
$labels = @("count", "mean", "std", "min", "25%", "50%", "75%", "max")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Match the bold/centered/bordered header formatting used on row 1 for the
# new index label column.
$ws.Range("B1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)

